# Update Name of Algo
# Applies corrected RandomForest imputation results to columns A and C
# for the rows flagged by the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "A4" = -20.4386
    "C4" = -12.0742
    "A6" = -22.58870000000002
    "A7" = -19.05389999999998
    "C9" = -10.2091
    "C12" = -10.46829999999999
    "A16" = -21.93130000000001
    "C17" = -14.37099999999998
    "C18" = -13.29309999999999
    "C19" = -11.2403
    "A20" = -20.10099999999998
    "C20" = -12.4884
    "C26" = -11.3272
    "A28" = -22.0573
    "A29" = -21.36659999999998
    "C31" = -12.7893
    "A32" = -21.216
    "C39" = -11.0885
    "A40" = -20.2587
    "C40" = -12.38060000000001
    "C41" = -12.08879999999999
    "C42" = -11.4148
    "C43" = -12.81179999999999
    "A46" = -21.65560000000002
    "C47" = -12.06729999999999
    "C48" = -11.6457
    "A51" = -21.65849999999999
    "A52" = -22.3394
    "A57" = -22.38430000000002
    "A59" = -22.30040000000001
    "A62" = -22.20240000000003
    "C63" = -10.2931
    "C64" = -10.25189999999999
    "A66" = -21.4912
    "A73" = -20.16419999999997
    "A74" = -21.71639999999998
    "C76" = -12.199
    "C81" = -13.3157
    "C89" = -14.49119999999999
    "A92" = -21.6094
    "C94" = -10.5557
    "A100" = -21.98840000000001
}

foreach ($cell in $updates.Keys) {
    $ws.Range($cell).Value = $updates[$cell]
}
